# The presentation's "Simple Light" design theme is reverted back to the
# "Default" color theme (the deck's original design before "Simple Light"
# was applied). Only the 12-color theme palette changes; the font scheme
# and the fill/line/effect format scheme are identical between the two
# themes in this deck, so no other theme-level properties need updating.

function ConvertTo-RgbLong {
    # PowerPoint's RGBColor.RGB is a Long encoded as 0x00BBGGRR (i.e. the
    # bytes of a "RRGGBB" hex string in reverse order). Converts a plain
    # "RRGGBB" hex color string into that Long value.
    param([string]$HexColor)
    $r = [convert]::ToInt32($HexColor.Substring(0, 2), 16)
    $g = [convert]::ToInt32($HexColor.Substring(2, 2), 16)
    $b = [convert]::ToInt32($HexColor.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Default" theme color scheme, in msoThemeColorIndex order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5 Accent1 .. 10 Accent6,
# 11 Hyperlink, 12 FollowedHyperlink.
$defaultThemeColors = @(
    "000000", # Dark1
    "FFFFFF", # Light1
    "158158", # Dark2
    "F3F3F3", # Light2
    "058DC7", # Accent1
    "50B432", # Accent2
    "ED561B", # Accent3
    "EDEF00", # Accent4
    "24CBE5", # Accent5
    "64E572", # Accent6
    "2200CC", # Hyperlink
    "551A8B"  # FollowedHyperlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $defaultThemeColors.Length; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-RgbLong $defaultThemeColors[$i - 1]
}
